$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns H and I
$ws.Range("H1").Value = "Fecha Inicial Alquiler"
$ws.Range("I1").Value = "Fecha Final Alquiler"

# Row 2 (UML265 / Kia) - update Categoria text and add new Id Alquiler / date columns
$ws.Range("E2").Value = "pequeños Automóvil"
$ws.Range("F2").Value = "3"
$ws.Range("G2").Value = "1"
$ws.Range("H2").Value = "15/11/23"
$ws.Range("I2").Value = "16/11/23"

# Row 3 (LWJ823 / Kia / Rojo / Manual / SUV Automovil / 3 / 0 / 0/0/0 / 0/0/0)
$ws.Range("A3").Value = "LWJ823"
$ws.Range("B3").Value = "Kia"
$ws.Range("C3").Value = "Rojo"
$ws.Range("D3").Value = "Manual"
$ws.Range("E3").Value = "SUV Automóvil"
$ws.Range("F3").Value = "3"
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "0/0/0"
$ws.Range("I3").Value = "0/0/0"

# Row 4 (new): IKP84G / Pulsar / Gris / Manual / Naked Moto / 3 / 5 / 05/12/2023 / 07/12/2023
$ws.Range("A4").Value = "IKP84G"
$ws.Range("B4").Value = "Pulsar"
$ws.Range("C4").Value = "Gris"
$ws.Range("D4").Value = "Manual"
$ws.Range("E4").Value = "Naked Moto"
$ws.Range("F4").Value = "3"
$ws.Range("G4").Value = "5"
$ws.Range("H4").Value = "05/12/2023"
$ws.Range("I4").Value = "07/12/2023"

$ws.Range("J5").Select()
